$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-layout the header row ---
# Old layout: A S.No. | B Staff Name | C Mobile No. | D Date Of Birth | E Blood Group
#             F Designation | G RF ID Card No. | H Address Line - 1 | I Address Line - 2
#             J Address Line - 3 | K City | L Pincode | M Status
#
# New layout: A Staff Name | B Mobile No. | C Date Of Birth | D Blood Group | E Designation
#             F RF ID Card No. | G Address Line - 1 | H Address Line - 2 | I Address Line - 3
#             J City | K Pincode
# (the "S.No." and "Status" columns are removed)

$ws.Range("A1").Value = "Staff Name"
$ws.Range("B1").Value = "Mobile No."
$ws.Range("C1").Value = "Date Of Birth"
$ws.Range("D1").Value = "Blood Group"
$ws.Range("E1").Value = "Designation"
$ws.Range("F1").Value = "RF ID Card No."
$ws.Range("G1").Value = "Address Line - 1"
$ws.Range("H1").Value = "Address Line - 2"
$ws.Range("I1").Value = "Address Line - 3"
$ws.Range("J1").Value = "City"
$ws.Range("K1").Value = "Pincode"

# Remove the now-unused trailing columns (old L1 "Pincode" / M1 "Status")
$ws.Range("L1:M1").Clear() | Out-Null

# --- Row 2 sample/data row formatting ---
# F2 gets the same text-format style as the existing G2 cell
$ws.Range("F2").NumberFormat = $ws.Range("G2").NumberFormat

# --- Column width ---
# Column A (Staff Name) becomes wider (target stored width 17.85546875;
# 17.0 is the input value that lands closest to it given this runtime's
# internal column-width rounding)
$ws.Columns(1).ColumnWidth = 17

# --- Sheet view / selection ---
# Scroll back to show column A and select F2
$ws.Range("F2").Select()
